$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.086.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.307.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.644'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.19'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.64%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.654'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0990'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.77'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.11%  '
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.647.53'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.52'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.883'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.306.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.016.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000102'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.25%  '
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '239.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.29%  '
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.98%  '
$ws.Range("E28").Value = '  -1.63%  '
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0847'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.37%  '
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("E35").Value = '  +6.52%  '
$ws.Range("E36").Value = '  +1.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.86'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.36%  '
$ws.Range("E39").Value = '  -3.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '13.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.37'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.89'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.219'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.41%  '
$ws.Range("E44").Value = '  +2.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.93'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.38%  '
$ws.Range("E49").Value = '  -0.83%  '
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("E51").Value = '  -1.37%  '
